# Append three new bullet paragraphs (ListParagraph style, same numbered
# list as the existing "Notes about RIFF and WAV" bullets) after the last
# paragraph of the document ("Some WAVE files may have a byte alignment...").

$d = $word.ActiveDocument

# Locate the last paragraph in the document (the final existing bullet).
$lastPara = $d.Paragraphs.Last

# Create a new paragraph after it; Word carries over the ListParagraph
# style + numPr (numId 1 / ilvl 0) from the paragraph it was split from.
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$insertRange = $newPara.Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p1 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Should nChannels be less than the number of bits set in dwChannelMask, then the extra (most significant) bits in dwChannelMask are ignored (in WAVEFORMATEXTENSIBLE). Might be relevant, might not be. </w:t></w:r></w:p>'

$p2 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>WAVE_FORMAT_EXTENSIBLE is indeed a different format tag than WAVE_FORMAT_PCM, and it is required for sample rates above 16-bit from the sounds of it. I’ll support them both.</w:t></w:r></w:p>'

$p3 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>It looks like format chunks used to be different. I should make sure I’m reading the up-to-date format and maybe even supporting the old one too depending on how much it’s still in use</w:t></w:r><w:r><w:t xml:space="preserve"> (it isn’t).</w:t></w:r></w:p>'

$insertRange.InsertXML($p1 + $p2 + $p3)

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
